# Auto-generated PowerShell COM-interop script
# Updates cryptocurrency price/volume data in cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.359.79'
$ws.Range("E2").Value = '  -3.27%  '
# Row 3
$ws.Range("D3").Value = '2.243.03'
$ws.Range("E3").Value = '  -4.53%  '
# Row 4
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
$ws.Range("D5").Value = "'232.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.11%  '
# Row 6
$ws.Range("D6").Value = "'0.638"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.03%  '
# Row 7
$ws.Range("D7").Value = "'71.25"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.39%  '
# Row 8
$ws.Range("E8").Value = '  +0.07%  '
# Row 9
$ws.Range("E9").Value = '  -5.10%  '
# Row 10
$ws.Range("D10").Value = "'0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.65%  '
# Row 11
$ws.Range("D11").Value = "'58.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.86%  '
# Row 12
$ws.Range("D12").Value = "'35.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.41%  '
# Row 13
$ws.Range("E13").Value = '  -2.88%  '
# Row 14
$ws.Range("D14").Value = "'6.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.55%  '
# Row 15
$ws.Range("D15").Value = '2.576.68'
$ws.Range("E15").Value = '  -4.50%  '
# Row 16
$ws.Range("D16").Value = "'15.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.16%  '
# Row 17
$ws.Range("D17").Value = "'0.872"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.38%  '
# Row 18
$ws.Range("D18").Value = '2.253.41'
$ws.Range("E18").Value = '  -4.04%  '
# Row 19
$ws.Range("D19").Value = '42.177.28'
$ws.Range("E19").Value = '  -3.50%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("E20").Value = '  -2.99%  '
# Row 21
$ws.Range("D21").Value = "'73.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.09%  '
# Row 22
$ws.Range("D22").Value = "'6.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.67%  '
# Row 23
$ws.Range("D23").Value = "'237.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.20%  '
# Row 24
$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.45%  '
# Row 25
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.04%  '
# Row 26
$ws.Range("D26").Value = "'3.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.03%  '
# Row 27
$ws.Range("E27").Value = '  -5.07%  '
# Row 28
$ws.Range("D28").Value = "'10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.57%  '
# Row 29
$ws.Range("D29").Value = "'2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.92%  '
# Row 30
$ws.Range("D30").Value = "'167.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.19%  '
# Row 31
$ws.Range("D31").Value = "'20.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.57%  '
# Row 32
$ws.Range("E32").Value = '  -6.94%  '
# Row 33
$ws.Range("E33").Value = '  -6.37%  '
# Row 34
$ws.Range("D34").Value = "'5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.40%  '
# Row 35
$ws.Range("D35").Value = "'0.0722"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.22%  '
# Row 36
$ws.Range("E36").Value = '  -6.78%  '
# Row 37
$ws.Range("D37").Value = "'3.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.24%  '
# Row 38
$ws.Range("D38").Value = "'22.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +15.78%  '
# Row 39
$ws.Range("D39").Value = "'6.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.50%  '
# Row 40
$ws.Range("E40").Value = '  -4.96%  '
# Row 41
$ws.Range("D41").Value = "'0.0267"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.95%  '
# Row 42
$ws.Range("D42").Value = "'66.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.79%  '
# Row 43
$ws.Range("D43").Value = "'5.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.74%  '
# Row 44
$ws.Range("D44").Value = "'8.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.32%  '
# Row 45
$ws.Range("D45").Value = "'0.102"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.62%  '
# Row 46
$ws.Range("E46").Value = '  -5.75%  '
# Row 47
$ws.Range("E47").Value = '  +0.16%  '
# Row 48
$ws.Range("D48").Value = "'4.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.94%  '
# Row 49
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").Value = "'1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.88%  '
# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'2.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.50%  '
# Row 51
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").Value = "'10.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.53%  '
